# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (so tab order
#    becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2021-Q3, 2021-Q2, 2021-Q1,
#    2020-Q4) and populate it with the fund-holdings table for that
#    quarter.
# 2) Insert a new row 2 into the "总计" summary sheet for the "2022-Q4"
#    totals, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create + position the new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Match the look of the other quarter sheets: bold/bordered header row
# and a bold/bordered index column (column A), reusing the exact style
# already used by the "2022-Q3" sheet instead of inventing a new one.
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2").Copy()
$newSheet.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Header row ----
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# ---- Data rows 2..11 (index / 基金代码 / 基金名称 / 基金规模 /
#      股票总仓位 / 仓位占比 / 持有市值(亿元) / 仓位排名) ----
# Row 2
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'011431"
$newSheet.Cells.Item(2,3).Value = "'泰达宏利消费服务混合A"
$newSheet.Cells.Item(2,4).Value = "'2.00"
$newSheet.Cells.Item(2,5).Value = "'87.31"
$newSheet.Cells.Item(2,6).Value = "'3.45"
$newSheet.Cells.Item(2,7).Value = "'0.0690"
$newSheet.Cells.Item(2,8).Value = 6

# Row 3
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'011432"
$newSheet.Cells.Item(3,3).Value = "'泰达宏利消费服务混合C"
$newSheet.Cells.Item(3,4).Value = "'1.39"
$newSheet.Cells.Item(3,5).Value = "'87.31"
$newSheet.Cells.Item(3,6).Value = "'3.45"
$newSheet.Cells.Item(3,7).Value = "'0.0480"
$newSheet.Cells.Item(3,8).Value = 6

# Row 4
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'015784"
$newSheet.Cells.Item(4,3).Value = "'中信建投中证1000指数增强A"
$newSheet.Cells.Item(4,4).Value = "'6.76"
$newSheet.Cells.Item(4,5).Value = "'89.78"
$newSheet.Cells.Item(4,6).Value = "'0.58"
$newSheet.Cells.Item(4,7).Value = "'0.0392"
$newSheet.Cells.Item(4,8).Value = 9

# Row 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'002210"
$newSheet.Cells.Item(5,3).Value = "'创金合信量化多因子股票A"
$newSheet.Cells.Item(5,4).Value = "'3.02"
$newSheet.Cells.Item(5,5).Value = "'93.22"
$newSheet.Cells.Item(5,6).Value = "'1.28"
$newSheet.Cells.Item(5,7).Value = "'0.0387"
$newSheet.Cells.Item(5,8).Value = 1

# Row 6
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'013466"
$newSheet.Cells.Item(6,3).Value = "'博时智选量化多因子股票C"
$newSheet.Cells.Item(6,4).Value = "'2.26"
$newSheet.Cells.Item(6,5).Value = "'93.55"
$newSheet.Cells.Item(6,6).Value = "'1.36"
$newSheet.Cells.Item(6,7).Value = "'0.0307"
$newSheet.Cells.Item(6,8).Value = 7

# Row 7
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'004194"
$newSheet.Cells.Item(7,3).Value = "'招商中证1000指数增强A"
$newSheet.Cells.Item(7,4).Value = "'2.57"
$newSheet.Cells.Item(7,5).Value = "'94.27"
$newSheet.Cells.Item(7,6).Value = "'1.05"
$newSheet.Cells.Item(7,7).Value = "'0.0270"
$newSheet.Cells.Item(7,8).Value = 8

# Row 8
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'004195"
$newSheet.Cells.Item(8,3).Value = "'招商中证1000指数增强C"
$newSheet.Cells.Item(8,4).Value = "'2.14"
$newSheet.Cells.Item(8,5).Value = "'94.27"
$newSheet.Cells.Item(8,6).Value = "'1.05"
$newSheet.Cells.Item(8,7).Value = "'0.0225"
$newSheet.Cells.Item(8,8).Value = 8

# Row 9
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'003865"
$newSheet.Cells.Item(9,3).Value = "'创金合信量化多因子股票C"
$newSheet.Cells.Item(9,4).Value = "'1.68"
$newSheet.Cells.Item(9,5).Value = "'93.22"
$newSheet.Cells.Item(9,6).Value = "'1.28"
$newSheet.Cells.Item(9,7).Value = "'0.0215"
$newSheet.Cells.Item(9,8).Value = 1

# Row 10
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'015785"
$newSheet.Cells.Item(10,3).Value = "'中信建投中证1000指数增强C"
$newSheet.Cells.Item(10,4).Value = "'2.40"
$newSheet.Cells.Item(10,5).Value = "'89.78"
$newSheet.Cells.Item(10,6).Value = "'0.58"
$newSheet.Cells.Item(10,7).Value = "'0.0139"
$newSheet.Cells.Item(10,8).Value = 9

# Row 11
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "'013465"
$newSheet.Cells.Item(11,3).Value = "'博时智选量化多因子股票A"
$newSheet.Cells.Item(11,4).Value = "'0.67"
$newSheet.Cells.Item(11,5).Value = "'93.55"
$newSheet.Cells.Item(11,6).Value = "'1.36"
$newSheet.Cells.Item(11,7).Value = "'0.0091"
$newSheet.Cells.Item(11,8).Value = 7

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q4" row into the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert()
$totalSheet.Rows("2:2").ClearFormats()

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 10
$totalSheet.Cells.Item(2,4).Value = 0.32

# Re-apply the bold/bordered "index column" style to A2 (it was cleared
# above along with the rest of the inserted row) by reusing the style
# already present a few rows down, then restore A2's own value.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$totalSheet.Cells.Item(2,1).Value = 0
